$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D='62.930.66'; DNeedsText=$False; E='  -2.20%  '},
    @{Row=3; D='2.681.06'; DNeedsText=$False; E='  -2.87%  '},
    @{Row=4; D=$null; DNeedsText=$False; E='  +0.00%  '},
    @{Row=5; D='548.61'; DNeedsText=$True; E='  -4.98%  '},
    @{Row=6; D='157.42'; DNeedsText=$True; E='  -1.80%  '},
    @{Row=7; D='0.999'; DNeedsText=$True; E='  +0.09%  '},
    @{Row=8; D='0.589'; DNeedsText=$True; E='  -2.29%  '},
    @{Row=9; D=$null; DNeedsText=$False; E='  -4.49%  '},
    @{Row=10; D=$null; DNeedsText=$False; E='  -2.58%  '},
    @{Row=11; D=$null; DNeedsText=$False; E='  -4.80%  '},
    @{Row=12; D='5.12'; DNeedsText=$True; E='  -12.73%  '},
    @{Row=13; D='3.156.79'; DNeedsText=$False; E='  -2.85%  '},
    @{Row=14; D='26.00'; DNeedsText=$True; E='  -4.97%  '},
    @{Row=15; D='62.794.67'; DNeedsText=$False; E='  -1.86%  '},
    @{Row=16; D=$null; DNeedsText=$False; E='  -3.93%  '},
    @{Row=17; D='2.684.04'; DNeedsText=$False; E='  -3.06%  '},
    @{Row=18; D='11.91'; DNeedsText=$True; E='  -2.13%  '},
    @{Row=19; D=$null; DNeedsText=$False; E='  -5.71%  '},
    @{Row=20; D='343.09'; DNeedsText=$True; E='  -4.29%  '},
    @{Row=21; D=$null; DNeedsText=$False; E='  -5.28%  '},
    @{Row=22; D=$null; DNeedsText=$False; E='  +0.01%  '},
    @{Row=23; D=$null; DNeedsText=$False; E='  -4.99%  '},
    @{Row=24; D='63.44'; DNeedsText=$True; E='  -2.51%  '},
    @{Row=25; D=$null; DNeedsText=$False; E='  -2.01%  '},
    @{Row=26; D=$null; DNeedsText=$False; E='  +0.16%  '},
    @{Row=27; D='8.14'; DNeedsText=$True; E='  -5.47%  '},
    @{Row=28; D='0.0₃0852'; DNeedsText=$False; E='  -7.56%  '},
    @{Row=29; D=$null; DNeedsText=$False; E='  -2.68%  '},
    @{Row=30; D=$null; DNeedsText=$False; E='  -3.37%  '},
    @{Row=31; D=$null; DNeedsText=$False; E='  -4.95%  '},
    @{Row=32; D='165.38'; DNeedsText=$True; E='  -1.83%  '},
    @{Row=33; D=$null; DNeedsText=$False; E='  +0.02%  '},
    @{Row=34; D=$null; DNeedsText=$False; E='  -3.99%  '},
    @{Row=35; D='19.52'; DNeedsText=$True; E='  -3.41%  '},
    @{Row=36; D=$null; DNeedsText=$False; E='  -6.05%  '},
    @{Row=37; D=$null; DNeedsText=$False; E='  -3.60%  '},
    @{Row=38; D='338.65'; DNeedsText=$True; E='  -3.63%  '},
    @{Row=39; D='6.17'; DNeedsText=$True; E='  -3.61%  '},
    @{Row=40; D=$null; DNeedsText=$False; E='  -7.41%  '},
    @{Row=41; D='38.12'; DNeedsText=$True; E='  -2.22%  '},
    @{Row=42; D=$null; DNeedsText=$False; E='  -6.24%  '},
    @{Row=43; D='20.32'; DNeedsText=$True; E='  -5.83%  '},
    @{Row=44; D='20.74'; DNeedsText=$True; E='  -7.87%  '},
    @{Row=45; D=$null; DNeedsText=$False; E='  -2.09%  '},
    @{Row=46; D=$null; DNeedsText=$False; E='  -5.94%  '},
    @{Row=47; D='0.999'; DNeedsText=$True; E='  +0.09%  '},
    @{Row=48; D=$null; DNeedsText=$False; E='  +0.01%  '},
    @{Row=49; D='0.0974'; DNeedsText=$True; E='  -3.74%  '},
    @{Row=50; D='129.33'; DNeedsText=$True; E='  -5.39%  '},
    @{Row=51; D='2.085.78'; DNeedsText=$False; E='  -2.94%  '}
)

foreach ($item in $data) {
    $eCell = $ws.Cells.Item($item.Row, 5)
    $eCell.Value = $item.E
    if ($null -ne $item.D) {
        $dCell = $ws.Cells.Item($item.Row, 4)
        if ($item.DNeedsText) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $item.D
    }
}
